$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename (English column names) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Normalize "de/del/la/las/los/el/y" connector capitalization in place names ---
$edits = @{
    "B6" = "Pabellón De Arteaga"
    "B7" = "Rincón De Romos"
    "B8" = "San Francisco De Los Romo"
    "B32" = "Amatenango De La Frontera"
    "B34" = "Bejucal De Ocampo"
    "B39" = "Chiapa De Corzo"
    "B42" = "Comitán De Domínguez"
    "B63" = "Ocozocoautla De Espinosa"
    "B69" = "San Cristóbal De Las Casas"
    "B103" = "Coyame Del Sotol"
    "B114" = "Guadalupe Y Calvo"
    "B117" = "Hidalgo Del Parral"
    "B138" = "San Francisco De Borja"
    "B139" = "San Francisco Del Oro"
    "B146" = "Valle De Zaragoza"
    "B175" = "Villa De Álvarez"
    "A177" = "Ciudad De México"
    "B181" = "Cuajimalpa De Morelos"
    "B196" = "Coneto De Comonfort"
    "B210" = "Nombre De Dios"
    "B214" = "Pánuco De Coronado"
    "B221" = "San Juan Del Río"
    "B222" = "San Luis Del Cordero"
    "B223" = "San Pedro Del Gallo"
    "A233" = "Estado De México"
    "B233" = "Acambay De Ruíz Castañeda"
    "B236" = "Almoloya De Alquisiras"
    "B237" = "Almoloya De Juárez"
    "B243" = "Atizapán De Zaragoza"
    "B247" = "Chapa De Mota"
    "B251" = "Coacalco De Berriozábal"
    "B256" = "Ecatepec De Morelos"
    "B260" = "Ixtapan De La Sal"
    "B261" = "Ixtapan Del Oro"
    "B270" = "Naucalpan De Juárez"
    "B278" = "San Felipe Del Progreso"
    "B280" = "Soyaniquilpan De Juárez"
    "B289" = "Tenango Del Valle"
    "B296" = "Tlalnepantla De Baz"
    "B302" = "Valle De Bravo"
    "B303" = "Valle De Chalco Solidaridad"
    "B304" = "Villa De Allende"
    "B305" = "Villa Del Carbón"
    "B316" = "San Miguel De Allende"
    "B317" = "Apaseo El Alto"
    "B325" = "Dolores Hidalgo Cuna De La Independencia Nacional"
    "B329" = "Jaral Del Progreso"
    "B337" = "Purísima Del Rincón"
    "B341" = "San Diego De La Unión"
    "B343" = "San Francisco Del Rincón"
    "B345" = "San Luis De La Paz"
    "B346" = "Santa Cruz De Juventino Rosas"
    "B347" = "Silao De La Victoria"
    "B352" = "Valle De Santiago"
    "B357" = "Acapulco De Juárez"
    "B359" = "Ajuchitlán Del Progreso"
    "B362" = "Atenango Del Río"
    "B363" = "Atlamajalcingo Del Monte"
    "B365" = "Atoyac De Álvarez"
    "B366" = "Ayutla De Los Libres"
    "B369" = "Buenavista De Cuéllar"
    "B370" = "Chilapa De Álvarez"
    "B371" = "Chilpancingo De Los Bravo"
    "B372" = "Coahuayutla De José María Izazaga"
    "B377" = "Coyuca De Benítez"
    "B378" = "Coyuca De Catalán"
    "B382" = "Cuetzala Del Progreso"
    "B383" = "Cutzamala De Pinzón"
    "B388" = "Huitzuco De Los Figueroa"
    "B389" = "Iguala De La Independencia"
    "B391" = "Ixcateopan De Cuauhtémoc"
    "B394" = "La Unión De Isidoro Montes De Oca"
    "B397" = "Mártir De Cuilapan"
    "B410" = "Taxco De Alarcón"
    "B412" = "Técpan De Galeana"
    "B414" = "Tepecoacuilco De Trujano"
    "B416" = "Tixtla De Guerrero"
    "B419" = "Tlalixtaquilla De Maldonado"
    "B420" = "Tlapa De Comonfort"
    "B430" = "Agua Blanca De Iturbide"
    "B436" = "Atotonilco El Grande"
    "B440" = "Cuautepec De Hinojosa"
    "B444" = "Huasca De Ocampo"
    "B446" = "Huejutla De Reyes"
    "B449" = "Jacala De Ledezma"
    "B452" = "Mineral De La Reforma"
    "B453" = "Mineral Del Chico"
    "B454" = "Mineral Del Monte"
    "B455" = "Mixquiahuala De Juárez"
    "B456" = "Omitlán De Juárez"
    "B457" = "Pachuca De Soto"
    "B459" = "Progreso De Obregón"
    "B463" = "Santiago De Anaya"
    "B464" = "Santiago Tulantepec De Lugo Guerrero"
    "B468" = "Tenango De Doria"
    "B470" = "Tepehuacán De Guerrero"
    "B471" = "Tepeji Del Río De Ocampo"
    "B473" = "Tezontepec De Aldama"
    "B477" = "Tula De Allende"
    "B478" = "Tulancingo De Bravo"
    "B481" = "Zacualtipán De Ángeles"
    "B485" = "Acatlán De Juárez"
    "B486" = "Ahualulco De Mercado"
    "B489" = "Atemajac De Brizuela"
    "B492" = "Atotonilco El Alto"
    "B494" = "Autlán De Navarro"
    "B506" = "Cuautitlán De García Barragán"
    "B514" = "Encarnación De Díaz"
    "B519" = "Huejuquilla El Alto"
    "B520" = "Ixtlahuacán De Los Membrillos"
    "B521" = "Ixtlahuacán Del Río"
    "B524" = "Jilotlán De Los Dolores"
    "B529" = "La Manzanilla De La Paz"
    "B530" = "Lagos De Moreno"
    "B542" = "San Cristóbal De La Barranca"
    "B543" = "San Diego De Alejandría"
    "B545" = "San Juan De Los Lagos"
    "B546" = "San Juanito De Escobedo"
    "B548" = "San Martín De Bolaños"
    "B550" = "San Miguel El Alto"
    "B551" = "San Sebastián Del Oeste"
    "B552" = "Santa María De Los Ángeles"
    "B555" = "Talpa De Allende"
    "B556" = "Tamazula De Gordiano"
    "B562" = "Teocuitatlán De Corona"
    "B563" = "Tepatitlán De Morelos"
    "B565" = "Tizapán El Alto"
    "B566" = "Tlajomulco De Zúñiga"
    "B575" = "Unión De San Antonio"
    "B576" = "Unión De Tula"
    "B577" = "Valle De Guadalupe"
    "B582" = "Yahualica De González Gallo"
    "B583" = "Zacoalco De Torres"
    "B586" = "Zapotitlán De Vadillo"
    "B587" = "Zapotlán El Grande"
    "B611" = "Coalcomán De Vázquez Pallares"
    "B613" = "Cojumatlán De Régules"
    "B676" = "Tiquicheo De Nicolás Romero"
    "B700" = "Coatlán Del Río"
    "B709" = "Puente De Ixtla"
    "B714" = "Tetela Del Volcán"
    "B716" = "Tlaltizapán De Zapata"
    "B726" = "Amatlán De Cañas"
    "B729" = "Ixtlán Del Río"
    "B736" = "Santa María Del Oro"
    "B752" = "San Nicolás De Los Garza"
    "B756" = "Acatlán De Pérez Figueroa"
    "B758" = "Ayoquezco De Aldama"
    "B762" = "Chalcatongo De Hidalgo"
    "B764" = "Coicoyán De Las Flores"
    "B765" = "Cuilápam De Guerrero"
    "B767" = "Guevea De Humboldt"
    "B768" = "Heroica Ciudad De Ejutla De Crespo"
    "B769" = "Heroica Ciudad De Huajuapan De León"
    "B770" = "Heroica Ciudad De Tlaxiaco"
    "B772" = "Ixtlán De Juárez"
    "B773" = "Heroica Ciudad De Juchitán De Zaragoza"
    "B777" = "Mariscala De Juárez"
    "B779" = "Mazatlán Villa De Flores"
    "B780" = "Miahuatlán De Porfirio Díaz"
    "B782" = "Oaxaca De Juárez"
    "B783" = "Ocotlán De Morelos"
    "B784" = "Putla Villa De Guerrero"
    "B788" = "San Agustín De Las Juntas"
    "B797" = "San Antonio De La Cal"
    "B803" = "San Felipe Jalapa De Díaz"
    "B811" = "San José Del Progreso"
    "B846" = "San Pedro El Alto"
    "B856" = "Santa Ana Del Valle"
    "B901" = "Tlacolula De Matamoros"
    "B902" = "Totontepec Villa De Morelos"
    "B903" = "Villa De Tamazulápam Del Progreso"
    "B904" = "Villa De Tututepec De Melchor Ocampo"
    "B905" = "Villa De Zaachila"
    "B907" = "Villa Sola De Vega"
    "B908" = "Villa Tejúpam De La Unión"
    "B911" = "Zimatlán De Álvarez"
    "B924" = "Ayotoxco De Guerrero"
    "B927" = "Chalchicomula De Sesma"
    "B940" = "Cuayuca De Andrade"
    "B946" = "Huehuetlán El Chico"
    "B949" = "Huitzilan De Serdán"
    "B952" = "Izúcar De Matamoros"
    "B961" = "Los Reyes De Juárez"
    "B965" = "Palmar De Bravo"
    "B975" = "San Salvador El Seco"
    "B976" = "San Salvador El Verde"
    "B980" = "Tecali De Herrera"
    "B985" = "Tepanco De López"
    "B988" = "Tepexi De Rodríguez"
    "B993" = "Tlacotepec De Benito Juárez"
    "B1013" = "Amealco De Bonfil"
    "B1015" = "Cadereyta De Montes"
    "B1019" = "Jalpan De Serra"
    "B1021" = "Pinal De Amoles"
    "B1024" = "San Juan Del Río"
    "B1031" = "Armadillo De Los Infante"
    "B1039" = "Mexquitic De Carmona"
    "B1043" = "San Ciro De Acosta"
    "B1046" = "Santa María Del Río"
    "B1048" = "Soledad De Graciano Sánchez"
    "B1051" = "Villa De Ramos"
    "B1052" = "Villa De Reyes"
    "B1115" = "Nacozari De García"
    "B1130" = "San Miguel De Horcasitas"
    "B1131" = "San Pedro De La Cueva"
    "B1170" = "Contla De Juan Cuamatzi"
    "B1173" = "Mazatecochco De José María Morelos"
    "B1176" = "Papalotla De Xicohténcatl"
    "B1183" = "Ziltlaltépec De Trinidad Sánchez Santos"
    "B1190" = "Amatlán De Los Reyes"
    "B1196" = "Castillo De Teayo"
    "B1198" = "Cazones De Herrera"
    "B1202" = "Chinampa De Gorostiza"
    "B1211" = "Cosamaloapan De Carpio"
    "B1220" = "Hueyapan De Ocampo"
    "B1221" = "Ignacio De La Llave"
    "B1224" = "Ixhuatlán De Madero"
    "B1225" = "Ixhuatlán Del Café"
    "B1231" = "Juchique De Ferrer"
    "B1234" = "Lerdo De Tejada"
    "B1236" = "Martínez De La Torre"
    "B1239" = "Medellín De Bravo"
    "B1249" = "Ozuluama De Mascareñas"
    "B1253" = "Paso De Ovejas"
    "B1257" = "Poza Rica De Hidalgo"
    "B1264" = "Sayula De Alemán"
    "B1305" = "Cañitas De Felipe Pescador"
    "B1308" = "El Plateado De Joaquín Amaro"
    "B1317" = "Jiménez Del Teul"
    "B1321" = "Mezquital Del Oro"
    "B1325" = "Moyahua De Estrada"
    "B1326" = "Nochistlán De Mejía"
    "B1327" = "Noria De Ángeles"
    "B1338" = "Teúl De González Ortega"
    "B1339" = "Tlaltenango De Sánchez Román"
    "B1342" = "Villa De Cos"
}
foreach ($cell in $edits.Keys) {
    $ws.Range($cell).Value = $edits[$cell]
}

# --- Re-sync the few percentage cells whose last binary digit shifted on recalculation ---
$pctFix = @{
    "D270" = 0.0009418242411070056
    "D309" = 0.0009418242411070056
    "D329" = 0.0009418242411070056
    "D361" = 0.0009418242411070056
    "D554" = 0.0009418242411070056
    "D563" = 0.0009418242411070056
    "D644" = 0.0009418242411070056
    "D733" = 0.0009418242411070056
    "D804" = 0.0009418242411070056
    "D864" = 0.0009418242411070056
    "D886" = 0.0009418242411070056
    "D1022" = 0.0009418242411070056
    "D1033" = 0.0009418242411070056
    "D1056" = 0.0009418242411070056
    "D1084" = 0.0009418242411070056
    "D1134" = 0.0009418242411070056
    "D1326" = 0.0009418242411070056
}
foreach ($cell in $pctFix.Keys) {
    $ws.Range($cell).Value = $pctFix[$cell]
}

# --- Remove trailing footnote rows (1349-1353); dimension recalculates to A1:D1347 ---
$ws.Rows("1349:1353").Delete()

Write-Output ("Edits applied: {0}, pct fixups: {1}" -f $edits.Count, $pctFix.Count)
